$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:F for rows 2-25
$arr1 = New-Object 'object[,]' 24,5
$arr1[0,0] = 1.02
$arr1[0,1] = 1.027283669358748
$arr1[0,2] = 1.036106916628989
$arr1[0,3] = 1.027402826265482
$arr1[0,4] = 1.044741476468274
$arr1[1,0] = 1.02
$arr1[1,1] = 1.028041358513458
$arr1[1,2] = 1.036706316148222
$arr1[1,3] = 1.028041325041443
$arr1[1,4] = 1.04557350670045
$arr1[2,0] = 1.02
$arr1[2,1] = 1.028532359969291
$arr1[2,2] = 1.037094819278008
$arr1[2,3] = 1.028455492218789
$arr1[2,4] = 1.046113102094098
$arr1[3,0] = 1.02
$arr1[3,1] = 1.028738948837004
$arr1[3,2] = 1.037258300040131
$arr1[3,3] = 1.028629849290427
$arr1[3,4] = 1.046340236962582
$arr1[4,0] = 1.02
$arr1[4,1] = 1.028773646038603
$arr1[4,2] = 1.037285758156379
$arr1[4,3] = 1.028659138695606
$arr1[4,4] = 1.046378390787668
$arr1[5,0] = 1.02
$arr1[5,1] = 1.028535119746992
$arr1[5,2] = 1.037097003113095
$arr1[5,3] = 1.028457821041305
$arr1[5,4] = 1.046116135948648
$arr1[6,0] = 1.02
$arr1[6,1] = 1.027539582492267
$arr1[6,2] = 1.03630935042813
$arr1[6,3] = 1.027618398398981
$arr1[6,4] = 1.045022411838935
$arr1[7,0] = 1.02
$arr1[7,1] = 1.025790963613829
$arr1[7,2] = 1.034926481627421
$arr1[7,3] = 1.026147095949343
$arr1[7,4] = 1.043104541981318
$arr1[8,0] = 1.02
$arr1[8,1] = 1.024629128480081
$arr1[8,2] = 1.034008097817604
$arr1[8,3] = 1.025171636296171
$arr1[8,4] = 1.041832421477441
$arr1[9,0] = 1.02
$arr1[9,1] = 1.024126993078187
$arr1[9,2] = 1.033611289080211
$arr1[9,3] = 1.02475055963014
$arr1[9,4] = 1.041283139230409
$arr1[10,0] = 1.02
$arr1[10,1] = 1.023940621874928
$arr1[10,2] = 1.033464027488403
$arr1[10,3] = 1.024594350957842
$arr1[10,4] = 1.04107934714268
$arr1[11,0] = 1.02
$arr1[11,1] = 1.023980592571516
$arr1[11,2] = 1.033495609646669
$arr1[11,3] = 1.024627849257493
$arr1[11,4] = 1.041123050530506
$arr1[12,0] = 1.02
$arr1[12,1] = 1.024111584622398
$arr1[12,2] = 1.033599113710366
$arr1[12,3] = 1.024737643322493
$arr1[12,4] = 1.041266288889548
$arr1[13,0] = 1.02
$arr1[13,1] = 1.024192312321167
$arr1[13,2] = 1.033662903377845
$arr1[13,3] = 1.024805317361682
$arr1[13,4] = 1.041354574083279
$arr1[14,0] = 1.02
$arr1[14,1] = 1.024662473671987
$arr1[14,2] = 1.034034450945071
$arr1[14,3] = 1.025199609413967
$arr1[14,4] = 1.041868908473287
$arr1[15,0] = 1.02
$arr1[15,1] = 1.024957648372688
$arr1[15,2] = 1.034267743918025
$arr1[15,3] = 1.025447288927368
$arr1[15,4] = 1.042191954591433
$arr1[16,0] = 1.02
$arr1[16,1] = 1.025129910008451
$arr1[16,2] = 1.034403902365097
$arr1[16,3] = 1.025591881856504
$arr1[16,4] = 1.042380531743955
$arr1[17,0] = 1.02
$arr1[17,1] = 1.02518866218789
$arr1[17,2] = 1.034450342826606
$arr1[17,3] = 1.02564120555974
$arr1[17,4] = 1.042444857022069
$arr1[18,0] = 1.02
$arr1[18,1] = 1.024925969470487
$arr1[18,2] = 1.03424270523139
$arr1[18,3] = 1.025420702247144
$arr1[18,4] = 1.042157279287523
$arr1[19,0] = 1.02
$arr1[19,1] = 1.024073006730556
$arr1[19,2] = 1.033568630715399
$arr1[19,3] = 1.024705306225918
$arr1[19,4] = 1.041224102235546
$arr1[20,0] = 1.02
$arr1[20,1] = 1.02353755072862
$arr1[20,2] = 1.033145571312085
$arr1[20,3] = 1.024256654651356
$arr1[20,4] = 1.040638742151998
$arr1[21,0] = 1.02
$arr1[21,1] = 1.023821326224974
$arr1[21,2] = 1.03336977058669
$arr1[21,3] = 1.024494383994376
$arr1[21,4] = 1.040948922426933
$arr1[22,0] = 1.02
$arr1[22,1] = 1.024940283521178
$arr1[22,2] = 1.034254018881148
$arr1[22,3] = 1.025432715235546
$arr1[22,4] = 1.042172947102175
$arr1[23,0] = 1.02
$arr1[23,1] = 1.026242342663042
$arr1[23,2] = 1.035283372879811
$arr1[23,3] = 1.026526518286922
$arr1[23,4] = 1.043599228421567
$ws.Range("B2:F25").Value = $arr1

# Columns I:N for rows 2-25
$arr2 = New-Object 'object[,]' 24,6
$arr2[0,0] = 1.034647283600106
$arr2[0,1] = 1.032442402237036
$arr2[0,2] = 1.038901605891513
$arr2[0,3] = 1.030222657885666
$arr2[0,4] = 1.047511673474632
$arr2[0,5] = 1.014840840667054
$arr2[1,0] = 1.034795896174332
$arr2[1,1] = 1.032841196576826
$arr2[1,2] = 1.039311096567165
$arr2[1,3] = 1.030669332004703
$arr2[1,4] = 1.048154949072271
$arr2[1,5] = 1.014973815848604
$arr2[2,0] = 1.03489121677391
$arr2[2,1] = 1.03309928336749
$arr2[2,2] = 1.039576018439458
$arr2[2,3] = 1.030958687131356
$arr2[2,4] = 1.048571757974476
$arr2[2,5] = 1.015059845897843
$arr2[3,0] = 1.034931087309368
$arr2[3,1] = 1.033207791358892
$arr2[3,2] = 1.039687379187725
$arr2[3,3] = 1.0310804088764
$arr2[3,4] = 1.048747118215047
$arr2[3,5] = 1.01509600915785
$arr2[4,0] = 1.034937769866322
$arr2[4,1] = 1.033226010776817
$arr2[4,2] = 1.039706076379114
$arr2[4,3] = 1.031100850957039
$arr2[4,4] = 1.048776569761754
$arr2[4,5] = 1.01510208089335
$arr2[5,0] = 1.034891750321569
$arr2[5,1] = 1.033100733225603
$arr2[5,2] = 1.039577506497734
$arr2[5,3] = 1.030960313283139
$arr2[5,4] = 1.048574100623273
$arr2[5,5] = 1.015060329128512
$arr2[6,0] = 1.034697681740138
$arr2[6,1] = 1.03257716753924
$arr2[6,2] = 1.039040003860586
$arr2[6,3] = 1.030373544889444
$arr2[6,4] = 1.047728953041839
$arr2[6,5] = 1.014885782822518
$arr2[7,0] = 1.034349294966848
$arr2[7,1] = 1.031654943898597
$arr2[7,2] = 1.038092565341473
$arr2[7,3] = 1.029342153761343
$arr2[7,4] = 1.046244116683139
$arr2[7,5] = 1.014578125242933
$arr2[8,0] = 1.034112770061868
$arr2[8,1] = 1.031040452970536
$arr2[8,2] = 1.03746082835244
$arr2[8,3] = 1.028656374763775
$arr2[8,4] = 1.045257310110844
$arr2[8,5] = 1.014372991508891
$arr2[9,0] = 1.034009349564539
$arr2[9,1] = 1.030774465358734
$arr2[9,2] = 1.037187270632216
$arr2[9,3] = 1.028359873932318
$arr2[9,4] = 1.044830767174712
$arr2[9,5] = 1.0142841655249
$arr2[10,0] = 1.0339707845571
$arr2[10,1] = 1.030675680745035
$arr2[10,2] = 1.037085658837726
$arr2[10,3] = 1.028249808837239
$arr2[10,4] = 1.044672444862262
$arr2[10,5] = 1.014251171826388
$arr2[11,0] = 1.033979063659102
$arr2[11,1] = 1.030696869678202
$arr2[11,2] = 1.037107454901984
$arr2[11,3] = 1.028273415049086
$arr2[11,4] = 1.044706400313258
$arr2[11,5] = 1.014258249068656
$arr2[12,0] = 1.034006164825173
$arr2[12,1] = 1.030766299477346
$arr2[12,2] = 1.037178871366908
$arr2[12,3] = 1.028350774510812
$arr2[12,4] = 1.044817677851484
$arr2[12,5] = 1.014281438248053
$arr2[13,0] = 1.034022842890302
$arr2[13,1] = 1.030809079499978
$arr2[13,2] = 1.037222873414361
$arr2[13,3] = 1.028398447352421
$arr2[13,4] = 1.044886254875215
$arr2[13,5] = 1.014295725911407
$arr2[14,0] = 1.03411961265849
$arr2[14,1] = 1.031058107725541
$arr2[14,2] = 1.037478983359486
$arr2[14,3] = 1.028676062073017
$arr2[14,4] = 1.045285634343871
$arr2[14,5] = 1.014378886602083
$arr2[15,0] = 1.034180045605437
$arr2[15,1] = 1.031214341875931
$arr2[15,2] = 1.037639632317697
$arr2[15,3] = 1.028850322925485
$arr2[15,4] = 1.04553635690729
$arr2[15,5] = 1.014431050988027
$arr2[16,0] = 1.034215198232139
$arr2[16,1] = 1.031305479303431
$arr2[16,2] = 1.037733334937141
$arr2[16,3] = 1.028952009263911
$arr2[16,4] = 1.045682671451073
$arr2[16,5] = 1.01446147739296
$arr2[17,0] = 1.034227167912371
$arr2[17,1] = 1.031336556243907
$arr2[17,2] = 1.037765284861501
$arr2[17,3] = 1.028986688913484
$arr2[17,4] = 1.045732573136657
$arr2[17,5] = 1.014471851962236
$arr2[18,0] = 1.03417357173837
$arr2[18,1] = 1.03119757852148
$arr2[18,2] = 1.037622396318574
$arr2[18,3] = 1.028831621933725
$arr2[18,4] = 1.045509449255365
$arr2[18,5] = 1.014425454260158
$arr2[19,0] = 1.033998188342713
$arr2[19,1] = 1.030745853706259
$arr2[19,2] = 1.037157840994076
$arr2[19,3] = 1.028327992173998
$arr2[19,4] = 1.044784906207556
$arr2[19,5] = 1.014274609602533
$arr2[20,0] = 1.033887050189997
$arr2[20,1] = 1.0304619240056
$arr2[20,2] = 1.036865755733148
$arr2[20,3] = 1.028011737309557
$arr2[20,4] = 1.044330021956434
$arr2[20,5] = 1.014179769132158
$arr2[21,0] = 1.03394604861751
$arr2[21,1] = 1.0306124317547
$arr2[21,2] = 1.037020595259867
$arr2[21,3] = 1.028179351847964
$arr2[21,4] = 1.044571101054761
$arr2[21,5] = 1.014230045550694
$arr2[22,0] = 1.034176497299705
$arr2[22,1] = 1.031205153133115
$arr2[22,2] = 1.037630184528564
$arr2[22,3] = 1.02884007197448
$arr2[22,4] = 1.045521607441861
$arr2[22,5] = 1.014427983181253
$arr2[23,0] = 1.03444011668045
$arr2[23,1] = 1.031893309325373
$arr2[23,2] = 1.038337526661241
$arr2[23,3] = 1.029608479556313
$arr2[23,4] = 1.046627446998191
$arr2[23,5] = 1.014657669223838
$ws.Range("I2:N25").Value = $arr2
